# Add a new "2014-07-15" sensor log sheet, placed right after "Main",
# as a snapshot/duplicate of the current "Main" sheet's data
# (Building Groups / sensor layout), mirroring how the previous
# "2014-07-08" sheet was created from "Main" the week before.

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("Main")

# Copy "Main" and drop the copy immediately after it.
$mainSheet.Copy($null, $mainSheet)

# Excel names the new copy "Main (2)"; rename it to the new date-stamped name.
$newSheet = $wb.Worksheets.Item("Main (2)")
$newSheet.Name = "2014-07-15"

# Keep "Main" as the selected/active sheet, same as before the edit.
$mainSheet.Activate()
